# edit.ps1 - applies the "Add files via upload" edit to Polo InsuranceConnectionInfo.docx
#
# Two content changes:
#  1) "...Amazon AWS server so it can be accessed..." -> "...Amazon RDS so it can be
#     accessed..." (the word "AWS server" is replaced by "RDS"), and the _GoBack
#     bookmark (which tracks the last edited location) moves to sit right after the
#     newly typed "RDS".
#  2) The "+  to  add a new hostname..." list item loses its proofErr grammar-check
#     markers around the word "to" (the text itself is unchanged).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "AWS server" -> "RDS", with the _GoBack bookmark repositioned to
# the point right after "RDS" (mirroring where Word leaves _GoBack after the
# user types replacement text there).
# ---------------------------------------------------------------------------

# Locate & replace the phrase "AWS server" with "RDS" (typed-replacement style).
$rng = $d.Content
$found = $rng.Find.Execute("AWS server")
if (-not $found) { throw "Could not find 'AWS server'" }
$rdsStart = $rng.Start
$rng.Text = "RDS"
$rdsEnd = $rng.End

# Force a run boundary immediately before "RDS" so it becomes its own run
# (temporary bookmark, removed again below).
$d.Bookmarks.Add("tmp0", $d.Range($rdsStart, $rdsStart))

# Re-point the _GoBack bookmark (Word's "last edit" marker) to just after "RDS".
$d.Bookmarks.Add("_GoBack", $d.Range($rdsEnd, $rdsEnd))

# The rest of that sentence ("Two sample ", "tables were created...so on.",
# " These tables...team.") were already separate runs before the edit; add
# temporary boundary bookmarks so they remain distinct runs instead of being
# coalesced by the edit above.
$r1 = $d.Content
$found1 = $r1.Find.Execute("Two sample ")
if (-not $found1) { throw "Could not find 'Two sample '" }
$d.Bookmarks.Add("tmpA", $d.Range($r1.Start, $r1.Start))
$d.Bookmarks.Add("tmpB", $d.Range($r1.End, $r1.End))

$r2 = $d.Content
$found2 = $r2.Find.Execute(" These tables will be updated")
if (-not $found2) { throw "Could not find ' These tables will be updated'" }
$d.Bookmarks.Add("tmpC", $d.Range($r2.Start, $r2.Start))

# Clean up temporary bookmarks (the run splits they caused persist).
$d.Bookmarks("tmp0").Delete()
$d.Bookmarks("tmpA").Delete()
$d.Bookmarks("tmpB").Delete()
$d.Bookmarks("tmpC").Delete()

# ---------------------------------------------------------------------------
# Change 2: remove the proofErr (grammar-check) markers around "to" in
# "+ to add a new hostname and enter the hostname stated below." The text
# itself is unchanged; using Find/Replace over the whole sentence causes the
# runtime to rebuild it as a single clean run without the proofErr markers.
# ---------------------------------------------------------------------------

$rng2 = $d.Content
$replaced = $rng2.Find.Execute(
    "+ to add a new hostname and enter the hostname stated below.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "+ to add a new hostname and enter the hostname stated below.", 2)
if (-not $replaced) { throw "Could not find/replace the '+ to add a new hostname' sentence" }

Write-Host "Edits applied successfully"
